$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "saveAs" link/value for the newly added column F, header row 1.
# (Excel stores ColumnWidth internally with a small offset vs. the displayed
#  "width" attribute in the XML; 48.1666... maps to a stored width of 49,
#  matching the target column width exactly.)
$ws.Range("F1").Value = "//?page=1&rows=2&newname=`$_otg&cols=3,4&val=5"
$ws.Range("F1").ColumnWidth = 48.166666666666664

# The active selection moves from the old F7 (outside data) to the newly
# added F1 cell, reflecting that saveAs is now only clickable after the
# "calculated" column has been filled in.
$ws.Range("F1").Select()
